# Update cryptocurrency price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.665.19'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.918.65'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("D4").Formula = '''1.002'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Formula = '''239.46'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").Formula = '''1.001'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Formula = '''0.4939'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Formula = '''0.2975'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Formula = '''0.06754'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").Value = '1.895.95'
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").Formula = '''17.14'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Formula = '''0.07343'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Formula = '''5.185'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("D14").Formula = '''88.69'
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("D15").Formula = '''0.6708'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '30.642.17'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Formula = '''0.000007948'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Formula = '''13.49'
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '2.164.66'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").Formula = '''5.359'
$ws.Range("E21").Value = '  +11.10%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Formula = '''200.70'
$ws.Range("E23").Value = '  +4.92%  '
$ws.Range("D24").Formula = '''6.314'
$ws.Range("E24").Value = '  +2.90%  '
$ws.Range("D25").Formula = '''9.631'
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("D26").Formula = '''164.86'
$ws.Range("E26").Value = '  +6.49%  '
$ws.Range("D27").Formula = '''18.68'
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("D28").Formula = '''1.962'
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("D29").Formula = '''1.481'
$ws.Range("E29").Value = '  +5.67%  '
$ws.Range("D30").Formula = '''4.378'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Formula = '''0.09181'
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D32").Formula = '''4.056'
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("D33").Formula = '''0.05282'
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Formula = '''0.7427'
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").Formula = '''1.115'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").Formula = '''2.723'
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").Formula = '''0.01846'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Formula = '''2.723'
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").Formula = '''0.9260'
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").Formula = '''2.084'
$ws.Range("E40").Value = '  -3.06%  '
$ws.Range("D41").Formula = '''0.4461'
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Formula = '''72.72'
$ws.Range("E42").Value = '  +25.76%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Formula = '''5.955'
$ws.Range("E43").Value = '  +3.07%  '
$ws.Range("D44").Formula = '''106.46'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").Formula = '''1.003'
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").Formula = '''0.1395'
$ws.Range("E46").Value = '  +3.36%  '
$ws.Range("D47").Formula = '''7.619'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Formula = '''9.072'
$ws.Range("E48").Value = '  +3.62%  '
$ws.Range("D49").Formula = '''35.22'
$ws.Range("E49").Value = '  +4.53%  '
$ws.Range("D50").Formula = '''0.05887'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Formula = '''0.4027'
$ws.Range("E51").Value = '  +2.12%  '
